# Applies the changes described by the commit:
#  - adds a new "Player Info" worksheet (placed before "ODI Batting")
#    containing one header row (ID, NAME, BATTING_HAND, BOWL_STYLE)
#    and one data row describing player 4824.
#  - renames the MATCH_CARD_LINK column on "ODI Batting" to MATCH_CODE
#    and replaces the full scorecard URLs with just the numeric match
#    code that used to be the query-string parameter.
#
# NOTE: worksheet variables captured via Worksheets.Item(<index>) behave
# as *positional* references in this runtime, so once a new sheet is
# inserted before "ODI Batting" its old index-1 handle would silently
# start pointing at the new sheet instead. To stay safe, every sheet is
# re-fetched by its (stable) *name* right before it is used.

$wb = $excel.ActiveWorkbook

$odiName = "ODI Batting"

# ---------------------------------------------------------------------
# 1. Insert a brand-new "Player Info" sheet ahead of "ODI Batting".
# ---------------------------------------------------------------------
$info = $wb.Worksheets.Add($wb.Worksheets.Item($odiName))
$info.Name = "Player Info"

# Header row - reuse the same look & feel (bold, bordered, centered)
# already used for the header row on "ODI Batting".
$info = $wb.Worksheets.Item("Player Info")
$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

$hdr = $info.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Data row describing player 4824. The ID is stored as text (matching
# the inline-string representation used throughout these sheets), so
# it is entered with a leading apostrophe to keep it from turning into
# a number.
$info.Range("A2").Value = "'4824"
$info.Range("B2").Value = "Benjamin Reginald Mcdermott"
$info.Range("C2").Value = "Right Handed"
$info.Range("D2").Value = "Right Arm Medium"

# ---------------------------------------------------------------------
# 3. Update "ODI Batting": rename column D header and replace the full
#    match-scorecard links with their bare numeric match codes.
# ---------------------------------------------------------------------
$odi = $wb.Worksheets.Item($odiName)
$odi.Range("D1").Value = "MATCH_CODE"

$odi.Range("D2").Value = "'4483"
$odi.Range("D3").Value = "'4484"
$odi.Range("D4").Value = "'4564"
$odi.Range("D5").Value = "'4565"
$odi.Range("D6").Value = "'4567"
